$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("turnou")
$ws.Name = "knnpp"

$ws.Range("B2").Value = "K óptimo knnpp"
$ws.Range("C2").Value = 15

$ws.Range("B3").Value = "Arreglo aleatorio óptimo knnpp"
$ws.Range("C3").Value = 73

$ws.Range("B4").Value = "MAE knnpp"
$ws.Range("C4").Value = 0.7066677097061078

$ws.Range("B5").Value = "MSE knnpp"
$ws.Range("C5").Value = 1.165834308865758

$ws.Range("B6").Value = "RMSE knnpp"
$ws.Range("C6").Value = 1.079738074194736

$ws.Range("B7").Value = "R-cuadrado knnpp"
$ws.Range("C7").Value = 0.1742914713936639

$ws.Range("B8:C10").Clear()
